# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" sheets, as published by the latest site build.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value for the "展览" sheet (sheetId 1)
$exhibitionUpdates = @{
    7  = 14496
    9  = 655
    10 = 15121
    11 = 27
    12 = 8537
    13 = 304
    15 = 54
    22 = 48
    23 = 1066
    33 = 251
    34 = 411
    36 = 5251
    37 = 5226
}

# Row -> new F-column value for the "全部类型" sheet (sheetId 4)
$allTypesUpdates = @{
    7  = 14496
    9  = 655
    10 = 15121
    11 = 27
    12 = 8537
    13 = 304
    16 = 54
    23 = 48
    24 = 1066
    36 = 251
    37 = 411
    39 = 5251
    40 = 5226
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
